$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total 70 -> 56, Wrong marking -1 -> -2
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -2

# Row 12 E: score string "70 / 140" -> "54 / 112"
$ws.Range("E12").Value = "54 / 112"
